$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the quantity of part "Resistance" (row 3, column C) from 1 to 3
$ws.Range("C3").Value = 3

# Move/restore the active selection to C4 (matches the author's last cursor position)
$ws.Range("C4").Select()
